$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-data")
$ws.Activate()

# Fix the example to refer to strains by their correct, zero-padded name
$ws.Range("A2").Value = "JJS-MGP001"
$ws.Range("A6").Value = "JJS-MGP001"
$ws.Range("A4").Value = "JJS-MGP020"

$ws.Range("B26").Select()
